# Implemented CfgFindKey routine. Config file is now parsed and written into
# the connection menu.
#
# This script:
#  1. Inserts a new worksheet named "Match" right after "CfgList".
#  2. Renames the existing "Sheet2" worksheet to "Compare".
#  3. Populates "Match" with a small 2x2 lookup table (de/hl headers with
#     their resolved values) using a bold header row and a text number
#     format so hex-like strings such as "E32F" are not reinterpreted.
#  4. Freezes the header row on "Match" and makes it the active sheet/tab,
#     matching the selection state Excel would leave behind after such an
#     edit (mirroring CfgList's own frozen header-row pane).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Match" sheet right after "CfgList" --------------
$cfgList = $wb.Worksheets.Item("CfgList")
$matchSheet = $wb.Worksheets.Add([System.Type]::Missing, $cfgList)
$matchSheet.Name = "Match"

# --- 2. Rename "Sheet2" to "Compare" -------------------------------------
$compareSheet = $wb.Worksheets.Item("Sheet2")
$compareSheet.Name = "Compare"

# --- 3. Fill in the Match sheet contents ---------------------------------
# Write column-by-column (A1, A2, then B1, B2) so new shared-string entries
# are interned in the same order the lookup values are produced: the key
# ("de"/"E32F" pair) before the resolved address ("hl"/"C0B0" pair).
$matchSheet.Range("A1").Value = "de"
$matchSheet.Range("A2").Value = "E32F"
$matchSheet.Range("B1").Value = "hl"
$matchSheet.Range("B2").Value = "C0B0"

# Bold header row, then apply a text number format to the value row first
# and the header row second so the style table ends up in the same order
# (bold/general, normal/text, bold/text).
$matchSheet.Rows("1:1").Font.Bold = $true
$matchSheet.Range("A2:B2").NumberFormat = "@"
$matchSheet.Range("A1:B1").NumberFormat = "@"

# --- 4. Freeze the header row, select B3, and activate the Match tab ----
$matchSheet.Range("A2").Select()
$matchSheet.Application.ActiveWindow.FreezePanes = $true
$matchSheet.Range("B3").Select()
